$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-13: Rv_ID, annot_int, Gene names, Function [CC], cluster_labels
$data = @(
    @("Rv2192c", 5, "trpD Rv2192c MTCY190.03c", "FUNCTION: Catalyzes the transfer of the phosphoribosyl group of 5-phosphorylribose-1-pyrophosphate (PRPP) to anthranilate to yield N-(5'-phosphoribosyl)-anthranilate (PRA).", 40),
    @("Rv1559", 4, "ilvA Rv1559 MTCY48.06c", "FUNCTION: Catalyzes the anaerobic formation of alpha-ketobutyrate and ammonia from threonine in a two-step reaction. The first step involved a dehydration of threonine and a production of enamine intermediates (aminocrotonate), which tautomerizes to its imine form (iminobutyrate). Both intermediates are unstable and short-lived. The second step is the nonenzymatic hydrolysis of the enamine/imine intermediates to form 2-ketobutyrate and free ammonia. In the low water environment of the cell, the second step is accelerated by RidA (By similarity). {ECO:0000250}.", 40),
    @("Rv1612", 3, "trpB Rv1612 MTCY01B2.04", "FUNCTION: The beta subunit is responsible for the synthesis of L-tryptophan from indole and L-serine. {ECO:0000250}.", 40),
    @("Rv0557", 3, "mgtA pimB Rv0557", "FUNCTION: Catalyzes the addition of a mannose residue from GDP-D-mannose to GlcAGroAc2 to generate 1,2-di-O-C16/C18:1-(alpha-D-mannopyranosyl)-(1-4)-(alpha-D-glucopyranosyluronic acid)-(1-3)-glycerol(ManGlcAGroAc2). {ECO:0000269|PubMed:19395496}.", 40),
    @("Rv1609", 3, "trpE Rv1609 MTCY01B2.01 MTV046.07", "FUNCTION: Part of a heterotetrameric complex that catalyzes the two-step biosynthesis of anthranilate, an intermediate in the biosynthesis of L-tryptophan. In the first step, the glutamine-binding beta subunit (TrpG) of anthranilate synthase (AS) provides the glutamine amidotransferase activity which generates ammonia as a substrate that, along with chorismate, is used in the second step, catalyzed by the large alpha subunit of AS (TrpE) to produce anthranilate. In the absence of TrpG, TrpE can synthesize anthranilate directly from chorismate and high concentrations of ammonia (By similarity). {ECO:0000250}.", 40),
    @("Rv3606c", 2, "folK Rv3606c MTCY07H7B.16", "", 40),
    @("Rv2367c", 2, "ybeY Rv2367c MTCY27.13", "FUNCTION: Single strand-specific metallo-endoribonuclease involved in late-stage 70S ribosome quality control and in maturation of the 3' terminus of the 16S rRNA. {ECO:0000255|HAMAP-Rule:MF_00009}.", 40),
    @("Rv1514c", 1, "Rv1514c MTCY277.36c", "", 40),
    @("Rv2645", 1, "Rv2645 MTCY441.15", "", 40),
    @("Rv1534", 1, "Rv1534", "", 40),
    @("Rv1879", 1, "Rv1879", "", 40),
    @("Rv1192", 1, "Rv1192", "", 40)
)

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $ws.Cells.Item($rowIndex, 4).Value = $entry[3]
    $ws.Cells.Item($rowIndex, 5).Value = $entry[4]
    $rowIndex++
}
